# Project Sample Project is saved.TEST Author: admin. Type: SAVE.
# Cell B11 on the active sheet ("Rules") changes from the text "R40"
# to the text "1" (entered as literal text, not a number).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B11").Value = "'1"
